$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$Address,
        [string]$Text
    )
    $rng = $ws.Range($Address)
    $rng.NumberFormat = "@"
    $rng.Value = $Text
    $rng.Style = "Normal"
}

# --- Simple price (column D) updates ---
Set-TextValue "D2"  "245.65"
Set-TextValue "D3"  "22.19"
Set-TextValue "D4"  "5.345"
Set-TextValue "D5"  "0.05903"
Set-TextValue "D6"  "3.395"
Set-TextValue "D8"  "0.8100"
Set-TextValue "D9"  "0.9635"
Set-TextValue "D11" "0.07387"
Set-TextValue "D12" "0.03425"
Set-TextValue "D13" "0.03028"
Set-TextValue "D14" "4.429"
Set-TextValue "D15" "0.09393"
Set-TextValue "D16" "0.001586"
Set-TextValue "D17" "0.04849"

# Row 18 volume label lost its "Worstin24h" suffix
Set-TextValue "E18" "17OneONE"

Set-TextValue "D19" "0.006089"
Set-TextValue "D20" "0.004080"
Set-TextValue "D21" "0.0009864"
Set-TextValue "D22" "0.00009704"
Set-TextValue "D23" "3.701"
Set-TextValue "D24" "2.165"

Set-TextValue "D40" "0.03928"

# --- Rows 41-43: coins rotated (Kick -> BKEX -> CEJI -> Kick) with refreshed prices/volumes ---
Set-TextValue "B41" "BKEXToken"
Set-TextValue "C41" "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D41" "0.1073"
Set-TextValue "E41" "40BKEXTokenBKK"

Set-TextValue "B42" "CEJI"
Set-TextValue "C42" "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D42" "0.002441"
Set-TextValue "E42" "41CEJICEJI"

Set-TextValue "B43" "KickToken"
Set-TextValue "C43" "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue "D43" "0.003040"
Set-TextValue "E43" "42KickTokenKICK"

Set-TextValue "D44" "0.005765"
Set-TextValue "D45" "0.00005301"

Set-TextValue "D48" "0.04992"
Set-TextValue "E48" "47BOLOBOLOWorstin24h"
